$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 3
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 7
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = 4
$ws.Range("G2").Value = 2
$ws.Range("H2").Value = 23
$ws.Range("I2").Value = 5

# Row 3
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 7
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 5
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 12
$ws.Range("I3").Value = 5

# Row 4
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = 8
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 5
$ws.Range("H4").Value = 56
$ws.Range("I4").Value = 5

# Row 5
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = 6
$ws.Range("E5").Value = 8
$ws.Range("F5").Value = 2
$ws.Range("G5").Value = 4
$ws.Range("H5").Value = 45
$ws.Range("I5").Value = 5

# Row 6 (new row)
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 5
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 3
$ws.Range("G6").Value = 3
$ws.Range("H6").Value = 34
$ws.Range("I6").Value = 5
$ws.Range("J6").Value = "train_dim2_1"

# Select I1 as the final active cell (matches sheetView selection in target)
$ws.Range("I1").Select()
